$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.045495732274246
$ws.Range("C2").Value = 1.511024478812717
$ws.Range("B3").Value = 7.446184443117565
$ws.Range("C3").Value = 2.966064903100272
$ws.Range("B4").Value = 8.365356156143278
$ws.Range("C4").Value = 4.676376492967628
$ws.Range("B5").Value = 10.8235484012382
$ws.Range("C5").Value = 5.974205067448837
$ws.Range("B6").Value = 11.10761367534823
$ws.Range("C6").Value = 7.267188880845531
$ws.Range("B7").Value = 11.4908013006829
$ws.Range("C7").Value = 8.89964018785696
$ws.Range("B8").Value = 13.19893784236318
$ws.Range("C8").Value = 10.32462859022819
$ws.Range("B9").Value = 13.78330691909713
$ws.Range("C9").Value = 11.73870016266258
$ws.Range("B10").Value = 14.63408189523207
$ws.Range("C10").Value = 13.18543951116577
$ws.Range("B11").Value = 15.17221522932397
$ws.Range("C11").Value = 14.64110091663916
$ws.Range("B12").Value = 16.16766562099805
$ws.Range("C12").Value = 15.81281664976013
$ws.Range("B13").Value = 23.88324367694838
$ws.Range("C13").Value = 17.25715788837302
$ws.Range("B14").Value = 29.5739133221779
$ws.Range("C14").Value = 18.74378175416534
$ws.Range("B15").Value = 31.77664218196445
$ws.Range("C15").Value = 20.0612437297103
$ws.Range("B16").Value = 33.21190198613429
$ws.Range("C16").Value = 21.45501244811165
$ws.Range("B17").Value = 38.27461502687395
$ws.Range("C17").Value = 22.84904724755951
$ws.Range("B18").Value = 40.551154937693
$ws.Range("C18").Value = 24.24146280720153
$ws.Range("B19").Value = 40.63672385303992
$ws.Range("C19").Value = 25.77376496883236
$ws.Range("B20").Value = 41.24145730389564
$ws.Range("C20").Value = 27.27051961136612
$ws.Range("B21").Value = 42.36673981791161
$ws.Range("C21").Value = 28.66138470800246
$ws.Range("B22").Value = 43.7443579894589
$ws.Range("C22").Value = 30.33734664837597
$ws.Range("B23").Value = 47.36880728150401
$ws.Range("C23").Value = 31.85728044171165
$ws.Range("B24").Value = 49.96512270353566
$ws.Range("C24").Value = 33.26448624559008
$ws.Range("B25").Value = 53.10863044606238
$ws.Range("C25").Value = 34.66600169719916
$ws.Range("B26").Value = 62.76763030291587
$ws.Range("C26").Value = 36.03638251444526
$ws.Range("B27").Value = 64.08613664236512
$ws.Range("C27").Value = 37.50806042685146
$ws.Range("B28").Value = 68.31734634790764
$ws.Range("C28").Value = 38.77124824876518
$ws.Range("B29").Value = 71.09117499775255
$ws.Range("C29").Value = 40.16745968938718
$ws.Range("B30").Value = 71.38599896212673
$ws.Range("C30").Value = 41.56823242021897
$ws.Range("B31").Value = 71.66442815299744
$ws.Range("C31").Value = 43.06685574205189
$ws.Range("B32").Value = 73.82119962265877
$ws.Range("C32").Value = 44.59725552689273
$ws.Range("B33").Value = 75.34405084818799
$ws.Range("C33").Value = 46.04582419601627
$ws.Range("B34").Value = 75.45257786781606
$ws.Range("C34").Value = 47.53016499544016
$ws.Range("B35").Value = 76.65167404255726
$ws.Range("C35").Value = 49.18294238678084
$ws.Range("B36").Value = 77.90682632833682
$ws.Range("C36").Value = 50.56594246333142
$ws.Range("B37").Value = 79.40600162663964
$ws.Range("C37").Value = 51.92835595680131
$ws.Range("B38").Value = 81.71255850454112
$ws.Range("C38").Value = 53.47260089445633
$ws.Range("B39").Value = 85.89332432902555
$ws.Range("C39").Value = 54.77535394211034
$ws.Range("B40").Value = 89.28421830128102
$ws.Range("C40").Value = 56.12861772990928
$ws.Range("B41").Value = 89.44248919341722
$ws.Range("C41").Value = 57.50301899865807
$ws.Range("B42").Value = 89.84359766055307
$ws.Range("C42").Value = 59.19899792347676
$ws.Range("B43").Value = 94.94724344014305
$ws.Range("C43").Value = 60.63303659559129
$ws.Range("B44").Value = 96.32508425569814
$ws.Range("C44").Value = 62.05762840974054
$ws.Range("B45").Value = 98.31700417647849
$ws.Range("C45").Value = 64.27234794557046
$ws.Range("B46").Value = 98.65919470832149
$ws.Range("C46").Value = 65.70392426160304
$ws.Range("B47").Value = 99.64160413208504
$ws.Range("C47").Value = 67.27432579685004

# Remove now-obsolete trailing rows 48:49
$ws.Range("A48:C49").Delete()
